$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $text) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Rows 12/13 (Solana / BinanceUSD) swap position in the list
$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextCell "D12" "21.34"
Set-TextCell "E12" "  +2.73%  "

$ws.Range("B13").Value = "BinanceUSD"
$ws.Range("C13").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextCell "D13" "0.9935"
Set-TextCell "E13" "  -0.64%  "

# Rows 20/21 (Avalanche / Dai) swap position in the list
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextCell "D20" "17.40"
Set-TextCell "E20" "  +3.43%  "

$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextCell "D21" "0.9964"
Set-TextCell "E21" "  -0.34%  "

# Plain price / volume refreshes for all remaining rows
Set-TextCell "D2" "28.905.70"
Set-TextCell "E2" "  +5.58%  "

Set-TextCell "D3" "1.810.36"
Set-TextCell "E3" "  +1.60%  "

Set-TextCell "D4" "0.9939"
Set-TextCell "E4" "  -0.62%  "

Set-TextCell "D5" "315.07"
Set-TextCell "E5" "  +0.28%  "

Set-TextCell "D6" "0.9946"
Set-TextCell "E6" "  -0.53%  "

Set-TextCell "D7" "0.5383"
Set-TextCell "E7" "  +2.51%  "

Set-TextCell "D8" "0.3822"
Set-TextCell "E8" "  +1.68%  "

Set-TextCell "D9" "0.07607"
Set-TextCell "E9" "  +2.58%  "

Set-TextCell "D10" "42.61"
Set-TextCell "E10" "  -0.20%  "

Set-TextCell "D11" "1.131"
Set-TextCell "E11" "  +3.20%  "

Set-TextCell "D14" "6.229"
Set-TextCell "E14" "  +1.66%  "

Set-TextCell "D15" "7.449"
Set-TextCell "E15" "  +6.46%  "

Set-TextCell "D16" "1.800.03"
Set-TextCell "E16" "  +1.19%  "

Set-TextCell "D17" "91.61"
Set-TextCell "E17" "  +2.11%  "

Set-TextCell "D18" "0.00001071"
Set-TextCell "E18" "  +1.26%  "

Set-TextCell "D19" "0.06433"
Set-TextCell "E19" "  -0.15%  "

Set-TextCell "D22" "5.989"
Set-TextCell "E22" "  +1.44%  "

Set-TextCell "D23" "28.834.66"
Set-TextCell "E23" "  +5.16%  "

Set-TextCell "D24" "11.44"
Set-TextCell "E24" "  +1.72%  "

Set-TextCell "E25" "  +1.79%  "

Set-TextCell "D26" "161.92"
Set-TextCell "E26" "  +3.91%  "

Set-TextCell "D27" "20.73"
Set-TextCell "E27" "  +2.56%  "

Set-TextCell "D28" "2.427"
Set-TextCell "E28" "  +2.61%  "

Set-TextCell "D29" "2.005.64"
Set-TextCell "E29" "  +0.96%  "

Set-TextCell "D30" "124.25"
Set-TextCell "E30" "  +2.36%  "

Set-TextCell "D31" "1.143"
Set-TextCell "E31" "  +4.52%  "

Set-TextCell "D32" "0.1019"
Set-TextCell "E32" "  +0.35%  "

Set-TextCell "D33" "5.807"
Set-TextCell "E33" "  +3.35%  "

Set-TextCell "D34" "3.663"
Set-TextCell "E34" "  +1.19%  "

Set-TextCell "D35" "0.2319"
Set-TextCell "E35" "  +12.73%  "

Set-TextCell "D36" "0.06588"
Set-TextCell "E36" "  +9.91%  "

Set-TextCell "D37" "0.02337"
Set-TextCell "E37" "  +3.40%  "

Set-TextCell "D38" "5.158"
Set-TextCell "E38" "  +5.22%  "

Set-TextCell "D39" "8.700"
Set-TextCell "E39" "  +5.94%  "

Set-TextCell "D40" "11.67"
Set-TextCell "E40" "  +2.81%  "

Set-TextCell "D41" "0.6385"
Set-TextCell "E41" "  +4.03%  "

Set-TextCell "D42" "1.233"
Set-TextCell "E42" "  +8.65%  "

Set-TextCell "D43" "0.9943"
Set-TextCell "E43" "  -0.51%  "

Set-TextCell "D44" "1.397"
Set-TextCell "E44" "  -2.64%  "

Set-TextCell "D45" "13.70"
Set-TextCell "E45" "  +3.89%  "

Set-TextCell "D46" "0.6001"
Set-TextCell "E46" "  +3.44%  "

Set-TextCell "D47" "3.682"
Set-TextCell "E47" "  +1.56%  "

Set-TextCell "D48" "126.19"
Set-TextCell "E48" "  +3.62%  "

Set-TextCell "D49" "2.007"
Set-TextCell "E49" "  +5.63%  "

Set-TextCell "D50" "1.164"
Set-TextCell "E50" "  +4.34%  "

Set-TextCell "D51" "0.06976"
Set-TextCell "E51" "  +3.59%  "
